# ---------------------------------------------------------------------------
# Homunkulus "database.xlsx" update
#   - Renames the first sheet "Tabelle1" -> "Backupplans"
#   - Adds a new "General" sheet after it, with a "Times Executed" counter
#   - Adds a second backup-plan row (G:\ -> long folder list) to Backupplans
#   - Restyles Backupplans (date col centered, label col centered, path col
#     word-wrapped) and widens column C / the sheet default column width
# ---------------------------------------------------------------------------

$xlCenter = -4108

$wb = $excel.ActiveWorkbook

# --- Sheet1: rename --------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Backupplans"

# --- Sheet2: create ----------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "General"

# --- Backupplans: new row of data ------------------------------------------
$ws1.Range("A2").Value = 44875
$ws1.Range("B2").Value = "G:\"
$paths = "C:\Users\Tim\Documents\.16151814`nC:\Users\Tim\Documents\.Autos`nC:\Users\Tim\Documents\.Bewerbungsunterlagen`nC:\Users\Tim\Documents\.Buchhaltung`nC:\Users\Tim\Documents\.Dounjinshi`nC:\Users\Tim\Documents\.Minecraft_Server`nC:\Users\Tim\Documents\.Projekte`nC:\Users\Tim\Documents\GitHub`nC:\Users\Tim\Pictures`nC:\Users\Tim\Videos`nC:\Users\Tim\Music"
$ws1.Range("C2").Value = $paths

# --- Backupplans: number format + alignment (keep numFmtId=14 "short date") -
$ws1.Range("A1:A2").NumberFormat = "mm-dd-yy"
$ws1.Range("A1:A2").VerticalAlignment = $xlCenter
$ws1.Range("A1:A2").NumberFormat = "mm-dd-yy"

$ws1.Range("B1:B2").HorizontalAlignment = $xlCenter
$ws1.Range("B1:B2").VerticalAlignment = $xlCenter

$ws1.Range("C1:C2").WrapText = $true
$ws1.Range("C1:C2").VerticalAlignment = $xlCenter

$ws1.Rows("1").RowHeight = 45

# --- Backupplans: small formatting leftovers / column sizing ---------------
$ws1.Columns("C").ColumnWidth = 117.140625
$ws1.Range("A1").EntireColumn.AutoFit() | Out-Null

# a lightly-touched, essentially empty row further down (matches upstream
# worksheet noise left over from manual selection/formatting)
$ws1.Range("A15").VerticalAlignment = $xlCenter
$ws1.Range("B15").HorizontalAlignment = $xlCenter
$ws1.Range("B15").VerticalAlignment = $xlCenter

# --- General: "Times Executed" counter --------------------------------------
$ws2.Range("A2").Value = "Times Executed "
$ws2.Range("A2").Font.Bold = $true
$ws2.Range("B2").Value = 1
$ws2.Range("B2").HorizontalAlignment = $xlCenter

$ws2.Columns("A").ColumnWidth = 15.285156
$ws2.Columns("B").ColumnWidth = 11.425781

# --- Selections / active sheet ---------------------------------------------
$ws1.Range("C15").Select()
$ws2.Range("L21").Select()
